$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'278.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'6.59%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'12"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'27.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'3.32%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'12"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'4.827"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.64%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'12"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.06353"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.92%"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "'12"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'6.947"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'3.57%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "'12"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'3.385"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'6.76%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'12"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'0.8804"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'3.41%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'12"
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").Value = "'0.9554"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'4.46%"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'12"
$ws.Range("G9").Style = "Normal"
$ws.Range("D10").Value = "'0.1477"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'5.07%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'12"
$ws.Range("G10").Style = "Normal"
$ws.Range("D11").Value = "'0.05147"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.13%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'12"
$ws.Range("G11").Style = "Normal"
$ws.Range("D12").Value = "'0.07276"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'2.69%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'12"
$ws.Range("G12").Style = "Normal"
$ws.Range("D13").Value = "'0.03147"
$ws.Range("D13").Style = "Normal"
$ws.Range("G13").Value = "'12"
$ws.Range("G13").Style = "Normal"
$ws.Range("D14").Value = "'0.09070"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.23%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'12"
$ws.Range("G14").Style = "Normal"
$ws.Range("D15").Value = "'0.001558"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.66%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'12"
$ws.Range("G15").Style = "Normal"
$ws.Range("D16").Value = "'0.0006287"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.66%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'12"
$ws.Range("G16").Style = "Normal"
$ws.Range("D17").Value = "'0.006081"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.26%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'12"
$ws.Range("G17").Style = "Normal"
$ws.Range("D18").Value = "'3.454"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.17%"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'12"
$ws.Range("G18").Style = "Normal"
$ws.Range("E19").Value = "'5.93%"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").Value = "'12"
$ws.Range("G19").Style = "Normal"
$ws.Range("D20").Value = "'0.3155"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'2.70%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'12"
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = "'0.1340"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'3.10%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'12"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'3.872"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-5.90%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'12"
$ws.Range("G22").Style = "Normal"
$ws.Range("E23").Value = "'1.71%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'12"
$ws.Range("G23").Style = "Normal"
$ws.Range("D24").Value = "'0.001177"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.04%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'12"
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = "'0.004286"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'6.00%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").Value = "'12"
$ws.Range("G25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001202"
$ws.Range("D26").Style = "Normal"
$ws.Range("G26").Value = "'12"
$ws.Range("G26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001691"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'3.13%"
$ws.Range("E27").Style = "Normal"
$ws.Range("G27").Value = "'12"
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = "'12"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'12"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'12"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'12"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'12"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'12"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'12"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'12"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'12"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'12"
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").Value = "'12"
$ws.Range("G38").Style = "Normal"
$ws.Range("G39").Value = "'12"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.04077"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.28%"
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").Value = "'12"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.006715"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'60.93%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'12"
$ws.Range("G41").Style = "Normal"
$ws.Range("E42").Value = "'4.71%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").Value = "'12"
$ws.Range("G42").Style = "Normal"
$ws.Range("E43").Value = "'2.81%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'12"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.01255"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-4.51%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").Value = "'12"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005233"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'2.00%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").Value = "'12"
$ws.Range("G45").Style = "Normal"
$ws.Range("E46").Value = "'0.05%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").Value = "'12"
$ws.Range("G46").Style = "Normal"
$ws.Range("D47").Value = "'2.377"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'820.84%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").Value = "'12"
$ws.Range("G47").Style = "Normal"
$ws.Range("D48").Value = "'0.02252"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'6.14%"
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").Value = "'12"
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002103"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.05%"
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").Value = "'12"
$ws.Range("G49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002003"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.05%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").Value = "'12"
$ws.Range("G50").Style = "Normal"
$ws.Range("G51").Value = "'12"
$ws.Range("G51").Style = "Normal"
